# company_keyword.xlsx — add "Jio" to the Reliance keyword list (REL row)
# and move the selection cursor down to B23 (matches the scrolled view
# the author had when they saved the workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 holds the "REL" company / its keyword list in column B.
# Append ",Jio" to the existing comma-separated keyword string.
$cell = $ws.Range("B12")
$current = $cell.Value2
if ($current -notmatch ",Jio(,|$)") {
    $cell.Value = $current + ",Jio"
}

# Reposition the viewport/selection the way the author left it: scrolled
# so row 22 is at the top, and the active cell/selection on B23.
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("B23").Select()
